$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row correct-answer value (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row correct/total marks (B12): 57 -> 95
$ws.Range("B12").Value = 95

# Update the Max column text on the Total row (E12): "55/84" -> "95/140"
$ws.Range("E12").Value = "95/140"
